$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look numeric (e.g. "0.999", "3.30").
# Force text format on the whole price column up front so Excel does not
# reinterpret/truncate these values as numbers when we write them below.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.374.55"
$ws.Range("E2").Value = "  -0.70%  "
$ws.Range("D3").Value = "1.655.49"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "213.07"
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("D6").Value = "0.529"
$ws.Range("E6").Value = "  +3.38%  "
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").Value = "23.57"
$ws.Range("E8").Value = "  +0.80%  "
$ws.Range("D9").Value = "0.260"
$ws.Range("E9").Value = "  +0.41%  "
$ws.Range("D10").Value = "0.0614"
$ws.Range("E10").Value = "  -0.94%  "
$ws.Range("D11").Value = "0.0908"
$ws.Range("E11").Value = "  +3.75%  "
$ws.Range("D12").Value = "1.889.19"
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("D13").Value = "1.651.88"
$ws.Range("E13").Value = "  -0.55%  "
$ws.Range("E14").Value = "  -1.02%  "
$ws.Range("D15").Value = "0.570"
$ws.Range("E15").Value = "  +4.07%  "
$ws.Range("D16").Value = "65.52"
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("D17").Value = "27.378.73"
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("D18").Value = "230.87"
$ws.Range("E18").Value = "  -6.45%  "
$ws.Range("D20").Value = "7.45"
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("D21").Value = "0.999"
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("D22").Value = "4.36"
$ws.Range("E22").Value = "  -2.61%  "
$ws.Range("D23").Value = "9.29"
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("D24").Value = "2.03"
$ws.Range("E24").Value = "  +0.65%  "
$ws.Range("D25").Value = "147.07"
$ws.Range("E25").Value = "  +0.79%  "
$ws.Range("D26").Value = "7.07"
$ws.Range("E26").Value = "  -1.45%  "
$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").Value = "0.113"
$ws.Range("E27").Value = "  +1.53%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "15.79"
$ws.Range("E28").Value = "  -2.56%  "
$ws.Range("B29").Value = "BinanceUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("E30").Value = "  -0.45%  "
$ws.Range("D31").Value = "1.19"
$ws.Range("E31").Value = "  -4.30%  "
$ws.Range("D32").Value = "3.30"
$ws.Range("E32").Value = "  -0.95%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "3.13"
$ws.Range("E33").Value = "  +0.30%  "
$ws.Range("B34").Value = "Maker"
$ws.Range("C34").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D34").Value = "1.424.53"
$ws.Range("E34").Value = "  -1.26%  "
$ws.Range("E35").Value = "  +0.78%  "
$ws.Range("E36").Value = "  -0.31%  "
$ws.Range("D37").Value = "0.908"
$ws.Range("E37").Value = "  -2.40%  "
$ws.Range("D38").Value = "0.571"
$ws.Range("E38").Value = "  -1.23%  "
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("E40").Value = "  +1.06%  "
$ws.Range("D41").Value = "0.998"
$ws.Range("E41").Value = "  -0.24%  "
$ws.Range("D42").Value = "5.56"
$ws.Range("E42").Value = "  +2.96%  "
$ws.Range("D43").Value = "65.06"
$ws.Range("E43").Value = "  -5.68%  "
$ws.Range("E44").Value = "  +0.66%  "
$ws.Range("D45").Value = "0.790"
$ws.Range("E45").Value = "  -0.22%  "
$ws.Range("D46").Value = "1.797.64"
$ws.Range("E46").Value = "  -0.26%  "
$ws.Range("D47").Value = "1.67"
$ws.Range("E47").Value = "  -1.60%  "
$ws.Range("D48").Value = "88.11"
$ws.Range("E48").Value = "  -0.47%  "
$ws.Range("E49").Value = "  -2.28%  "
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("D51").Value = "7.77"
$ws.Range("E51").Value = "  -0.68%  "
